$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hit_miss_rule")

# Update the raw input values (columns H, I, J) which drive the
# ROUND()/SUM() formulas in columns D, E, F via recalculation.
$ws.Range("H5").Value = 88.698883056640625
$ws.Range("I6").Value = 11.30111694335938
$ws.Range("H8").Value = 5.2234883308410636
$ws.Range("I8").Value = 6.0832605361938477
$ws.Range("H9").Value = 8.42681884765625
$ws.Range("I9").Value = 9.2497806549072266
$ws.Range("H10").Value = 93.232131958007813
$ws.Range("I10").Value = 36.470077514648438
$ws.Range("J10").Value = 86.812469482421875

$excel.CalculateFullRebuild()
$wb.Save()
